$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.473.01"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "3.005.12"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'563.66"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'139.34"
$ws.Range("E6").Value = "  +3.99%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "2.992.85"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D11").Value = "'5.22"
$ws.Range("E11").Value = "  +7.22%  "
$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").Value = "'0.0000232"
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("D14").Value = "'33.83"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("E16").Value = "  +6.94%  "
$ws.Range("D17").Value = "3.498.27"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "3.000.42"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").Value = "59.366.25"
$ws.Range("E19").Value = "  +2.27%  "
$ws.Range("D20").Value = "'431.98"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").Value = "'13.66"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("E22").Value = "  +4.98%  "
$ws.Range("D23").Value = "'7.15"
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("D24").Value = "'13.52"
$ws.Range("E24").Value = "  +3.15%  "
$ws.Range("D25").Value = "'80.58"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'2.24"
$ws.Range("E27").Value = "  +10.94%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "'2.56"
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("D30").Value = "'7.92"
$ws.Range("E30").Value = "  +3.75%  "
$ws.Range("D31").Value = "'25.82"
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("D32").Value = "'6.15"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "'0.101"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  +5.73%  "
$ws.Range("E35").Value = "  +6.00%  "
$ws.Range("D36").Value = "0.0₃0762"
$ws.Range("E36").Value = "  +8.10%  "
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("D38").Value = "'49.00"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("D39").Value = "'8.69"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'412.99"
$ws.Range("E40").Value = "  +8.72%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.76"
$ws.Range("E41").Value = "  +6.05%  "
$ws.Range("D42").Value = "'0.0355"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "2.778.72"
$ws.Range("E43").Value = "  +3.30%  "
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").Value = "'0.254"
$ws.Range("E45").Value = "  +3.93%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'34.91"
$ws.Range("E47").Value = "  +20.64%  "
$ws.Range("D48").Value = "'123.78"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").Value = "'2.02"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "'23.61"
$ws.Range("E51").Value = "  -0.06%  "
